$wb = $excel.ActiveWorkbook

# "Generate Report for Archive" - refresh localization status: the handoff
# has moved on to translation, so update the Status cells and re-fit the
# Status columns to the new (shorter) text, on all three sheets.

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$newStatus = "In Translation"

$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsZhCn.Range("C2").Value = $newStatus
$wsDeDe.Range("C2").Value = $newStatus

# Target stored column width (from the xlsx diff) is ~13.41 chars; the
# ColumnWidth setter here snaps to the nearest 1/6-character pixel grid,
# so 12.5 is the closest input that lands on that value.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
